# Update "想去人数" (want-to-go count) values in column F for sheets
# "展览" and "全部类型" to reflect refreshed data as of commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F
$updates = @{
    3  = 2915
    7  = 1575
    11 = 1302
    13 = 420
    20 = 2961
    21 = 358
    23 = 76
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
